$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Recibido" (column B) dates for rows 11-17 -----------------------
# These rows only had a "Pedido" date (column A) before; the author now
# fills in the date the parts were actually received. Seed each B cell
# from its same-row A cell (same date style) and then overwrite with the
# real received date so the date number-format carries over.
$rows = 11..17
foreach ($r in $rows) {
    $ws.Range("A$r").Copy($ws.Range("B$r"))
    $ws.Range("B$r").Value = 40966
}

# --- Row 30: new line item - Hotend / extrusor -------------------------
$ws.Range("C30").Value = 1
$ws.Range("C30").Interior.ColorIndex = 6
$ws.Range("D30").Value = "Hotend"
$ws.Range("E30").Value = 50
$ws.Range("F30").Value = "http://www.hotends.com/index.php?route=product/product&product_id=70"

# --- Row 31: new line item - tabla/madera (wood board) ------------------
$ws.Range("A29").Copy($ws.Range("A31"))
$ws.Range("A31").Value = 41119
$ws.Range("A29").Copy($ws.Range("B31"))
$ws.Range("B31").Value = 40968

$ws.Range("C31").Value = 1
$ws.Range("C31").Interior.ColorIndex = 6

$ws.Range("F31").Value = "Maderas aparacio"

$ws.Range("C20").Copy($ws.Range("D31"))
$ws.Range("D31").Value = "Madera 22,5x22,5"

$ws.Range("C20").Copy($ws.Range("E31"))
$ws.Range("E31").Value = 0.5

# --- Restore the active cell selection left by the author --------------
$ws.Range("A29").Select()
